$d = $word.ActiveDocument

# Locate the unique occurrence of "imshow" in the sentence about displaying
# the image with Hough lines, and replace it with "imwrite".
$rng = $d.Content
$rng.Find.Execute("Hough lines superimposed using imshow", $false, $false, $false, $false, $false, $true, 1, $false, "Hough lines superimposed using imwrite", 2)
